$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 285.42856
$ws.Range("I38").Value = 19.8
$ws.Range("J38").Value = 949.5
$ws.Range("K38").Value = 59.40000000000001
$ws.Range("L38").Value = 2848.5
$ws.Range("M38").Value = 312.6
$ws.Range("N38").Value = -3592.5
$ws.Range("H43").Value = 2404.4285
$ws.Range("J43").Value = 3112.5
$ws.Range("L43").Value = 3112.5
$ws.Range("N43").Value = -3250.5
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = ""
$ws.Range("N48").Value = 0
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = ""
$ws.Range("N56").Value = 0
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = ""
$ws.Range("H70").Value = 3062.5
$ws.Range("I70").Value = 3062.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9187.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -8917.5
$ws.Range("H73").Value = 3062.5
$ws.Range("I73").Value = 3062.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9187.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -8251.5
$ws.Range("H113").Value = 5398.3335
$ws.Range("I113").Value = 4815.3335
$ws.Range("K113").Value = 4815.3335
$ws.Range("M113").Value = -1561.3335
$ws.Range("H116").Value = 10000
$ws.Range("I116").Value = 10000
$ws.Range("K116").Value = 10000
$ws.Range("M116").Value = -6558
$ws.Range("H129").Value = 2128
$ws.Range("J129").Value = 2299.3333
$ws.Range("L129").Value = 6897.999899999999
$ws.Range("N129").Value = -16897.9999
$ws.Range("H132").Value = 1474.25
$ws.Range("I132").Value = 1256.2858
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3768.8574
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1238.8574
$ws.Range("N132").Value = -14060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 183.85
$ws.Range("I5").Value = 183.26315
$ws.Range("K5").Value = 183.26315
$ws.Range("M5").Value = -71.26315
$ws.Range("H110").Value = 1197.6
$ws.Range("I110").Value = 747
$ws.Range("K110").Value = 747
$ws.Range("M110").Value = 1298

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 183.85
$ws.Range("I4").Value = 183.26315
$ws.Range("K4").Value = 183.26315
$ws.Range("M4").Value = -68.26315
$ws.Range("H86").Value = 5112.4287
$ws.Range("I86").Value = 5396.75
$ws.Range("J86").Value = 4733.3335
$ws.Range("K86").Value = 5396.75
$ws.Range("L86").Value = 4733.3335
$ws.Range("M86").Value = -4273.75
$ws.Range("N86").Value = -6979.3335
$ws.Range("H89").Value = 5112.4287
$ws.Range("I89").Value = 5396.75
$ws.Range("J89").Value = 4733.3335
$ws.Range("K89").Value = 26983.75
$ws.Range("L89").Value = 23666.6675
$ws.Range("M89").Value = -21367.75
$ws.Range("N89").Value = -34898.6675
$ws.Range("H107").Value = 600
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 600
$ws.Range("M107").Value = 1320

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 792.875
$ws.Range("I2").Value = 902
$ws.Range("J2").Value = 683.75
$ws.Range("K2").Value = 902
$ws.Range("L2").Value = 683.75
$ws.Range("M2").Value = -789
$ws.Range("N2").Value = -909.75
$ws.Range("H3").Value = 1749.75
$ws.Range("I3").Value = 1500
$ws.Range("J3").Value = 1999.5
$ws.Range("K3").Value = 1500
$ws.Range("L3").Value = 1999.5
$ws.Range("M3").Value = -1387
$ws.Range("N3").Value = -2225.5
$ws.Range("H5").Value = 1049.8334
$ws.Range("I5").Value = 699.75
$ws.Range("K5").Value = 699.75
$ws.Range("M5").Value = -587.75
$ws.Range("H7").Value = 103.46667
$ws.Range("I7").Value = 140.2
$ws.Range("K7").Value = 140.2
$ws.Range("M7").Value = -27.19999999999999
$ws.Range("H10").Value = 1064
$ws.Range("I10").Value = 1216.3334
$ws.Range("K10").Value = 1216.3334
$ws.Range("M10").Value = -1077.3334
$ws.Range("H12").Value = 169.75
$ws.Range("I12").Value = 169.75
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 169.75
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = ""
$ws.Range("N12").Value = 0.25
$ws.Range("H14").Value = 799
$ws.Range("J14").Value = 799
$ws.Range("L14").Value = 799
$ws.Range("N14").Value = -1139
$ws.Range("H15").Value = 9999.666999999999
$ws.Range("J15").Value = 9999.666999999999
$ws.Range("L15").Value = 9999.666999999999
$ws.Range("N15").Value = -10339.667
$ws.Range("H19").Value = 6000247
$ws.Range("I19").Value = 6000247
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 6000247
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = ""
$ws.Range("N19").Value = -6000077
$ws.Range("H22").Value = 990.75
$ws.Range("I22").Value = 999
$ws.Range("J22").Value = 949.5
$ws.Range("K22").Value = 999
$ws.Range("L22").Value = 949.5
$ws.Range("M22").Value = -649
$ws.Range("N22").Value = -1649.5
$ws.Range("H24").Value = 6000247
$ws.Range("I24").Value = 6000247
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 6000247
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = ""
$ws.Range("N24").Value = -6000077
$ws.Range("H62").Value = 4633.5
$ws.Range("J62").Value = 4600.3335
$ws.Range("L62").Value = 4600.3335
$ws.Range("N62").Value = -5848.3335
$ws.Range("H65").Value = 4633.5
$ws.Range("J65").Value = 4600.3335
$ws.Range("L65").Value = 23001.6675
$ws.Range("N65").Value = -29241.6675
$ws.Range("H132").Value = 2163.6667
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4989
$ws.Range("I80").Value = 4989
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 14967
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = ""
$ws.Range("N80").Value = -14031
$ws.Range("H83").Value = 4989
$ws.Range("I83").Value = 4989
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 44901
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = ""
$ws.Range("N83").Value = -40221
$ws.Range("H92").Value = 180.28572
$ws.Range("I92").Value = 145.33333
$ws.Range("K92").Value = 435.99999
$ws.Range("M92").Value = 812.00001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5499.5
$ws.Range("I70").Value = 5499.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 5499.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = ""
$ws.Range("N70").Value = -5229.5
$ws.Range("H73").Value = 5499.5
$ws.Range("I73").Value = 5499.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 5499.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = ""
$ws.Range("N73").Value = -4563.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 654.5
$ws.Range("I22").Value = 441.66666
$ws.Range("J22").Value = 782.2
$ws.Range("K22").Value = 441.66666
$ws.Range("L22").Value = 782.2
$ws.Range("M22").Value = -146.66666
$ws.Range("N22").Value = -1372.2
$ws.Range("H27").Value = 654.5
$ws.Range("I27").Value = 441.66666
$ws.Range("J27").Value = 782.2
$ws.Range("K27").Value = 441.66666
$ws.Range("L27").Value = 782.2
$ws.Range("M27").Value = -334.66666
$ws.Range("N27").Value = -996.2
$ws.Range("H32").Value = 1075.9
$ws.Range("I32").Value = 1075.9
$ws.Range("K32").Value = 1075.9
$ws.Range("M32").Value = -758.9000000000001
$ws.Range("H61").Value = 4397
$ws.Range("I61").Value = 3395
$ws.Range("J61").Value = 5399
$ws.Range("K61").Value = 3395
$ws.Range("L61").Value = 5399
$ws.Range("M61").Value = -3193
$ws.Range("N61").Value = -5803
$ws.Range("H68").Value = 3799.8
$ws.Range("I68").Value = 4124.75
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 4124.75
$ws.Range("L68").Value = 2500
$ws.Range("M68").Value = -3375.75
$ws.Range("N68").Value = -3998
$ws.Range("H71").Value = 3799.8
$ws.Range("I71").Value = 4124.75
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 20623.75
$ws.Range("L71").Value = 12500
$ws.Range("M71").Value = -16879.75
$ws.Range("N71").Value = -19988
$ws.Range("H113").Value = 4397
$ws.Range("I113").Value = 3395
$ws.Range("J113").Value = 5399
$ws.Range("K113").Value = 3395
$ws.Range("L113").Value = 5399
$ws.Range("M113").Value = -1225
$ws.Range("N113").Value = -9739
$ws.Range("H122").Value = 3238
$ws.Range("I122").Value = 2879.4119
$ws.Range("K122").Value = 8638.235700000001
$ws.Range("M122").Value = -6188.235700000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 527.0909
$ws.Range("I107").Value = 387.25
$ws.Range("K107").Value = 1161.75
$ws.Range("M107").Value = 758.25
